# Add "ContactUs" test-data sheet (Contact Us form test cases) after the
# existing sheets, matching the "Added contact us test cases" commit.

$wb = $excel.ActiveWorkbook

# Add the new worksheet as the LAST tab (after the current last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ContactUs"

# Header row + test-case rows for the Contact Us form.
$data = @(
    @("name",       "email",                  "subject",        "message",                        "expectedMessage"),
    @("John Doe",   "john@example.com",       "Test Subject 1", "This is a test message.",          "Success! Your details have been submitted successfully."),
    @("Jane Smith", "jane@example.com",       "Feedback",       "Loved the website!",                "Success! Your details have been submitted successfully."),
    @("Invalid",    "invalid-email",          "Test Error",     "Missing proper email format.",      "Please enter a valid email address."),
    @("Empty",      "",                       "No Email",       "No email provided.",                "Email is required."),
    @("Long Text",  "longtext@example.com",   "Test Long",      "Lorem ipsum dolor sit amet...",     "Success! Your details have been submitted successfully.")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        if ($value -ne "") {
            $ws.Cells.Item($r + 1, $c + 1).Value = $value
        }
    }
}

# Approximate the column widths Excel computed when the data was auto-fit.
$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 21.333333333333332
$ws.Columns.Item(5).ColumnWidth = 44.666666666666664

# Select the full used range, like the saved workbook shows.
$ws.Range("A1:E6").Select()
